# edit.ps1
# Applies the "Updated cryptos list" data refresh to the active worksheet.
# For every touched cell we force Text number format first so that
# numeric-looking strings (e.g. "302.19") are not silently converted
# into real numbers by Excel -- the source data model stores everything
# as plain text (inline strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.056.19'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.313.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.80%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.19'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.54'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.34%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.86%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.515'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.21'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +9.63%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.88'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +12.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.88'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.672.89'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.309.12'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.948.94'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.12%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +7.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.23'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.88'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.36'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.22'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +13.98%  '
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.71'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.96%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.76'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.82'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.02'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.60%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.42'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.01%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0695'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.23%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.29%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.94%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.08%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.985.94'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.46%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.69%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.22'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +7.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.92'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.63'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '56.14'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +8.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.539.23'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.01%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.02%  '
